$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.608.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.346.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.994"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.55%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.702.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.348.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.539.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "264.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.38%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0897"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  -9.00%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "120.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +22.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.75%  "
